# Sync automático del tracker (cada 3h)
# Appends the latest batch of tennis picks (2025-09-23) to the results
# tracker sheet, rows 61-66, and extends the used-range dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TrackerRow($Row, $EventId, $Fecha, $JugadorA, $JugadorB, $Pronostico, $Cuota) {
    # event_id / fecha look like numbers (or dates) to Excel's automatic
    # type inference, but the source feed stores them as plain text, so
    # force a text literal with a leading apostrophe and then strip the
    # resulting "quoted number" cell format back to the sheet's default
    # so no stray per-cell number format sticks around.
    $ws.Range("A$Row").Value = "'" + $EventId
    $ws.Range("A$Row").ClearFormats()

    $ws.Range("B$Row").Value = "'" + $Fecha
    $ws.Range("B$Row").ClearFormats()

    $ws.Range("C$Row").Value = $JugadorA
    $ws.Range("D$Row").Value = $JugadorB
    $ws.Range("E$Row").Value = $Pronostico
    $ws.Range("F$Row").Value = $Cuota

    # resultado / profit are still pending (the bet hasn't settled yet) -
    # write them as empty text cells, same as the rest of the open rows.
    $ws.Range("G$Row").Value = "'"
    $ws.Range("G$Row").ClearFormats()

    $ws.Range("H$Row").Value = "'"
    $ws.Range("H$Row").ClearFormats()
}

Add-TrackerRow 61 "14743054" "2025-09-23" "Ajla Tomljanovic"     "Yuliia Starodubtseva" "Gana Yuliia Starodubtseva" 2.2
Add-TrackerRow 62 "14743046" "2025-09-23" "Anastasija Sevastova" "Kimberly Birrell"     "Gana Kimberly Birrell"     2.2
Add-TrackerRow 63 "14743050" "2025-09-23" "Katerina Siniakova"   "Anastasia Potapova"   "Gana Anastasia Potapova"   2.5
Add-TrackerRow 64 "14743042" "2025-09-23" "Katie Boulter"        "Hailey Baptiste"      "Gana Hailey Baptiste"      1.91
Add-TrackerRow 65 "14743045" "2025-09-23" "Magda Linette"        "Bianca Andreescu"     "Gana Magda Linette"        2.1
Add-TrackerRow 66 "14743047" "2025-09-23" "Suzan Lamens"         "Yafan Wang"           "Gana Suzan Lamens"         1.53
